$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: populate column C (English label) for all new rows, top to bottom
# This establishes shared-string indices 907-919 in row order.
$ws.Range("C395").Value = 'Wood Heat; Existing'
$ws.Range("C396").Value = 'Wood and Electric Dual System; Existing'
$ws.Range("C397").Value = 'Wood and Oil Dual System; Existing'
$ws.Range("C398").Value = 'Furnace; Heating Oil; Existing'
$ws.Range("C399").Value = 'Oil Furnace'
$ws.Range("C400").Value = 'Furnace; Natural Gas; Existing'
$ws.Range("C401").Value = 'Natural Gas Furnace'
$ws.Range("C402").Value = 'Electric Baseboard; Existing'
$ws.Range("C403").Value = 'Ductless Air Source Heat Pump; Existing'
$ws.Range("C404").Value = 'Heat Pump (Ductless) with EBB backup'
$ws.Range("C405").Value = 'Heat Pump (Ductless) with oil backup'
$ws.Range("C406").Value = 'Heat Pump (Ductless) with wood backup'
$ws.Range("C407").Value = 'Heat Pump (Ductless) with NG backup'

# Step 2: populate column B (French label) only for rows 395-397, in the order 396, 395, 397
$ws.Range("B396").Value = 'Bi-énergie bois et électrique; Existant'
$ws.Range("B395").Value = 'Bois; Existant'
$ws.Range("B397").Value = 'Bi-énergie bois et mazout; Existant'

# Step 3: populate column A (variable) for rows 398-399, where it differs from column C
$ws.Range("A398").Value = 'Heating oil furnace; Existing'
$ws.Range("A399").Value = 'Heating oil furnace'

# Step 4: populate column A for the remaining rows, where A equals C (reuses existing shared string)
$ws.Range("A395").Value = 'Wood Heat; Existing'
$ws.Range("A396").Value = 'Wood and Electric Dual System; Existing'
$ws.Range("A397").Value = 'Wood and Oil Dual System; Existing'
$ws.Range("A400").Value = 'Furnace; Natural Gas; Existing'
$ws.Range("A401").Value = 'Natural Gas Furnace'
$ws.Range("A402").Value = 'Electric Baseboard; Existing'
$ws.Range("A403").Value = 'Ductless Air Source Heat Pump; Existing'
$ws.Range("A404").Value = 'Heat Pump (Ductless) with EBB backup'
$ws.Range("A405").Value = 'Heat Pump (Ductless) with oil backup'
$ws.Range("A406").Value = 'Heat Pump (Ductless) with wood backup'
$ws.Range("A407").Value = 'Heat Pump (Ductless) with NG backup'

# Step 5: populate column D (color code text) and copy the matching fill/format from a reference cell
$ws.Range("D395").Value = '#067906'
$ws.Range('D99').Copy()
$ws.Range("D395").PasteSpecial(-4122)
$ws.Range("D396").Value = '#13b157'
$ws.Range('D150').Copy()
$ws.Range("D396").PasteSpecial(-4122)
$ws.Range("D397").Value = '#5487a4'
$ws.Range('D2').Copy()
$ws.Range("D397").PasteSpecial(-4122)
$ws.Range("D398").Value = '#a4aba6'
$ws.Range('D166').Copy()
$ws.Range("D398").PasteSpecial(-4122)
$ws.Range("D399").Value = '#63625a'
$ws.Range('D69').Copy()
$ws.Range("D399").PasteSpecial(-4122)
$ws.Range("D400").Value = '#ffe36d'
$ws.Range('D121').Copy()
$ws.Range("D400").PasteSpecial(-4122)
$ws.Range("D401").Value = '#bfa904'
$ws.Range('D385').Copy()
$ws.Range("D401").PasteSpecial(-4122)
$ws.Range("D402").Value = '#9bd4dc'
$ws.Range('D89').Copy()
$ws.Range("D402").PasteSpecial(-4122)
$ws.Range("D403").Value = '#5487a4'
$ws.Range('D2').Copy()
$ws.Range("D403").PasteSpecial(-4122)
$ws.Range("D404").Value = '#314deb'
$ws.Range('D32').Copy()
$ws.Range("D404").PasteSpecial(-4122)
$ws.Range("D405").Value = '#f6b4a4'
$ws.Range('D17').Copy()
$ws.Range("D405").PasteSpecial(-4122)
$ws.Range("D406").Value = '#9bdb9a'
$ws.Range('D275').Copy()
$ws.Range("D406").PasteSpecial(-4122)
$ws.Range("D407").Value = '#bdd030'
$ws.Range('D24').Copy()
$ws.Range("D407").PasteSpecial(-4122)

$excel.CutCopyMode = $false
$ws.Range("D403").Select()
